# Generate Report for Handback
# Update the "Correspond Handback DateTime" (and refreshed handoff-cycle
# timestamp) for the second data row ("82b901f4-...") on both the zh-cn
# and de-de language sheets, reflecting a new handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 10:52:39"
$wsZhCn.Range("H3").Value = "2016-03-23 10:53:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 10:52:43"
$wsDeDe.Range("H3").Value = "2016-03-23 10:53:13"
